$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 1181
$ws1.Range("F4").Value = 16736
$ws1.Range("F5").Value = 28
$ws1.Range("F6").Value = 1636
$ws1.Range("F8").Value = 4
$ws1.Range("F10").Value = 212
$ws1.Range("F12").Value = 11616
$ws1.Range("F13").Value = 26
$ws1.Range("F14").Value = 1289
$ws1.Range("F15").Value = 4601
$ws1.Range("F16").Value = 427
$ws1.Range("F17").Value = 403
$ws1.Range("F20").Value = 334

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1
$ws4.Range("F3").Value = 1
$ws4.Range("F4").Value = 1181
$ws4.Range("F5").Value = 16736
$ws4.Range("F6").Value = 28
$ws4.Range("F7").Value = 1636
$ws4.Range("F9").Value = 4
$ws4.Range("F11").Value = 212
$ws4.Range("F15").Value = 11616
$ws4.Range("F16").Value = 26
$ws4.Range("F17").Value = 1289
$ws4.Range("F18").Value = 4601
$ws4.Range("F19").Value = 427
$ws4.Range("F20").Value = 403
$ws4.Range("F23").Value = 334
